$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (sheet1) ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 12763
$ws1.Range("F3").Value = 619
$ws1.Range("F5").Value = 28
$ws1.Range("F6").Value = 308
$ws1.Range("F7").Value = 400
$ws1.Range("F9").Value = 12758
$ws1.Range("F10").Value = 36
$ws1.Range("F11").Value = 8
$ws1.Range("F12").Value = 5201
$ws1.Range("F14").Value = 15
$ws1.Range("F16").Value = 25
$ws1.Range("F20").Value = 666
$ws1.Range("F22").Value = 6145
$ws1.Range("F23").Value = 1152
$ws1.Range("F24").Value = 3618
$ws1.Range("F26").Value = 39

# --- Sheet "全部类型" (sheet4) ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 12763
$ws4.Range("F3").Value = 619
$ws4.Range("F5").Value = 28
$ws4.Range("F6").Value = 308
$ws4.Range("F8").Value = 400
$ws4.Range("F10").Value = 12759
$ws4.Range("F11").Value = 36
$ws4.Range("F12").Value = 8
$ws4.Range("F13").Value = 5201
$ws4.Range("F15").Value = 15
$ws4.Range("F17").Value = 25
$ws4.Range("F21").Value = 666
$ws4.Range("F24").Value = 6145
$ws4.Range("F25").Value = 1152
$ws4.Range("F26").Value = 3618
$ws4.Range("F28").Value = 39
